$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $true, $false, $false, $false,
                       $true, 1, $false, $new, 2) | Out-Null
}

# 1. Intro paragraph: drop leading "Die " and change first "verbessern" -> "steigern"
Replace-Text "Die Fabrikam Inc. hat eine umfassende Initiative zur digitalen Transformation durchlaufen, die darauf abzielt, die betriebliche Effizienz zu verbessern, die Kundenerfahrung zu verbessern und Innovationen voranzutreiben." `
             "Fabrikam Inc. hat eine umfassende Initiative zur digitalen Transformation durchlaufen, die darauf abzielt, die betriebliche Effizienz zu steigern, die Kundenerfahrung zu verbessern und Innovationen voranzutreiben."

# 2. "Wichtige Updates" heading -> "Wichtige Aktualisierungen"
Replace-Text "Wichtige Updates" "Wichtige Aktualisierungen"

# 3. Fix "kigestützter" typo -> "KI-gesteuerter"
Replace-Text "Integration kigestützter Analysen zur Optimierung von Entscheidungsprozessen." `
             "Integration KI-gesteuerter Analysen zur Optimierung von Entscheidungsprozessen."

# 4. Chatbot bullet rewrite
Replace-Text "Einführung von Chatbots für den 24/7-Kundensupport, verringerung der Reaktionszeiten um 50%." `
             "Einführung von Chatbots für den 24/7-Kundensupport, wodurch die Reaktionszeiten um 50 % verkürzt werden."

# 5. Bold the "Prozessautomatisierung" heading
$rng = $d.Content
$found = $rng.Find.Execute("Prozessautomatisierung", $true, $true)
if ($found) {
    $rng.Font.Bold = 1
}

# 6. RPA bullet wording
Replace-Text "Implementierung der Roboterprozessautomatisierung (RPA) für Routineaufgaben." `
             "Implementierung der robotergesteuerten Prozessautomatisierung (RPA) für Routineaufgaben."

# 7. Processing time bullet rewrite
Replace-Text "Erreicht eine Reduzierung der Verarbeitungszeit für wichtige Geschäftsvorgänge um 40 %." `
             "Die Bearbeitungszeit für wichtige Geschäftsvorgänge konnte um 40 % reduziert werden."

# 8. "zu strategischeren Rollen" -> "auf strategischere Rollen"
Replace-Text "Neuverteilung von Personalressourcen zu strategischeren Rollen innerhalb der Organisation." `
             "Neuverteilung von Personalressourcen auf strategischere Rollen innerhalb der Organisation."

# 9. Gender-inclusive heading "Mitarbeitern" -> "Mitarbeiterinnen und Mitarbeitern"
Replace-Text "Schulung und Entwicklung von Mitarbeitern" "Schulung und Entwicklung von Mitarbeiterinnen und Mitarbeitern"

# 10. Gender-inclusive bullet "Mitarbeiter" -> "Mitarbeiterinnen und Mitarbeiter"
Replace-Text "Durchgeführte Digitalkompetenzprogramme für alle Mitarbeiter." `
             "Durchgeführte Digitalkompetenzprogramme für alle Mitarbeiterinnen und Mitarbeiter."

# 11. E-Learning bullet rewrite
Replace-Text "Eine neue E-Learning-Plattform mit Kursen zu neuen Technologien wurde gestartet." `
             "Start einer neuen E-Learning-Plattform mit Kursen zu neuen Technologien."

# 12. Mitarbeiterbindung bullet rewrite
Replace-Text "Erhöhte Mitarbeiterbindung und Einführung neuer Tools um 35 %." `
             "Verbesserung der Mitarbeitendenbindung und Einführung neuer Tools um 35 %."

# 13. Q1 milestone rewrite
Replace-Text "Q1 2024: Abgeschlossene Migration zur Cloudinfrastruktur." `
             "Q1 2024: Abschluss der Migration zur Cloudinfrastruktur."

# 14. Q2 milestone rewrite
Replace-Text "Q2 2024: Gestartete KI-gesteuerte Analyseplattform." `
             "Q2 2024: Start der KI-gesteuerten Analyseplattform."

# 15. Q3 milestone rewrite
Replace-Text "Q3 2024: Neue digitale Kundenportale eingeführt." `
             "Q3 2024: Einführung neuer digitaler Kundenportale."

# 16. Q4 milestone rewrite
Replace-Text "Q4 2024: Erreicht 50 % Automatisierung von Routineprozessen." `
             "Q4 2024: 50 % Automatisierung von Routineprozessen."

# 17. Bold the "Zukünftige Pläne" heading
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Zukünftige Pläne", $true, $true)
if ($found2) {
    $rng2.Font.Bold = 1
}

# 18. Future plans bullet: drop imperative "Sie" phrasing
Replace-Text "Erweitern Sie WEITERHIN KI- und Machine Learning-Anwendungen in allen Abteilungen." `
             "Kontinuierlicher Ausbau von KI- und Machine Learning-Anwendungen in allen Abteilungen."

Replace-Text "Verbessern Sie die digitale Kundenerfahrung mit neuen Features und Diensten weiter." `
             "Weitere Verbesserung der digitalen Kundenerfahrung mit neuen Features und Diensten."

Replace-Text "Konzentrieren Sie sich auf Cybersicherheitsmaßnahmen zum Schutz vor sich entwickelnden Bedrohungen." `
             "Konzentration auf Cybersicherheitsmaßnahmen zum Schutz vor sich entwickelnden Bedrohungen."

Replace-Text "Entwickeln Sie eine umfassende digitale Strategie für die nächsten fünf Jahre." `
             "Entwickeln einer umfassenden digitalen Strategie für die nächsten fünf Jahre."

# 19. Conclusion paragraph rewrite
Replace-Text "Die Digitale Transformationsreise der Fabrikam Inc. hat zu erheblichen Verbesserungen bei Effizienz, Kundenzufriedenheit und gesamter Geschäftsleistung geführt." `
             "Die digitale Transformation von Fabrikam Inc. hat zu erheblichen Verbesserungen bei der Effizienz, der Kundenzufriedenheit und der allgemeinen Unternehmensleistung geführt."

Replace-Text "Die Organisation setzt sich weiterhin dafür ein, Technologien zu nutzen, um zukünftiges Wachstum und Innovation voranzutreiben." `
             "Die Organisation setzt weiterhin auf die Nutzung von Technologien, um zukünftiges Wachstum und Innovation voranzutreiben."
